# Update "想去人数" (attendance count) values that changed between scrapes.
# Sheet "展览" (index 1): F3 294->296, F4 1291->1292, F6 61->62
# Sheet "全部类型" (index 4): F4 294->296, F5 1291->1292, F7 61->62

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 296
$wsExhibit.Range("F4").Value = 1292
$wsExhibit.Range("F6").Value = 62

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 296
$wsAll.Range("F5").Value = 1292
$wsAll.Range("F7").Value = 62
